$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Fitness column (C) values per contiguous Generation ranges,
# matching the commit's updated log values.
$ws.Range("C2:C11").Value = 8468
$ws.Range("C12:C16").Value = 8356
$ws.Range("C17:C28").Value = 8325
$ws.Range("C29:C31").Value = 7966
$ws.Range("C32:C98").Value = 7586
$ws.Range("C99:C252").Value = 7569
